# Auto-generated Word COM-interop script
# Rebuilds the resume body paragraph-by-paragraph to match the target OOXML.
$d = $word.ActiveDocument

# --- Step 1: clear all existing content down to a single clean paragraph ---
$existingCount = $d.Paragraphs.Count
$wipeRange = $d.Range(0, $d.Content.End)
for ($i = 0; $i -lt $existingCount; $i++) {
    $wipeRange.Delete()
}
$d.Paragraphs.Item(1).Style = "Normal"

# --- Step 2: paragraph data (text, style, alignment, bold, size-in-points) ---
$items = @(
    @{ Text = 'Dheeraj Chand'; Style = $null; Jc = 'center'; Bold = $true; Size = 14 },
    @{ Text = '202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/'; Style = $null; Jc = 'center'; Bold = $false; Size = $null },
    @{ Text = 'PROFESSIONAL SUMMARY'; Style = 'Heading 2'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Senior Software Engineer with 21 years building scalable geospatial data platforms, web applications, and distributed analytical systems. Expert in full-stack development with deep specialization in Apache Spark/Sedona for big data geospatial processing. Proven track record architecting multi-tenant SaaS platforms used by thousands of analysts, implementing ETL pipelines processing billions of geospatial records, and building production systems integrating ESRI, OSGeo, and SAFE FME technologies. Strong background in both enterprise consulting and startup environments, with experience leading engineering teams and delivering mission-critical geospatial applications.'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'CORE COMPETENCIES'; Style = 'Heading 2'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Programming & Development: Python: Django/GeoDjango, Flask, Pandas, PySpark, NumPy, SciKit-Learn • JVM: Scala (Spark/Sedona), Java (GeoTools, enterprise applications), Groovy • Web Technologies: JavaScript, React, d3.js, OpenLayers, jQuery, HTML/CSS • Database Languages: SQL, T-SQL, PostgreSQL/PostGIS, Oracle, MySQL • Statistical/Analysis: R, SPSS, NetLogo (agent-based modeling)'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Big Data & Geospatial Platforms: Apache Spark: PySpark, Spark SQL, Sedona (geospatial), distributed processing • Geospatial Stack: PostGIS, ESRI ArcGIS, Quantum GIS, GRASS, OSGeo, SAFE FME • Cloud Platforms: AWS (EC2, RDS, S3), Snowflake, Google Cloud, Microsoft Azure • Data Engineering: ETL/ELT pipelines, dbt, Hadoop, Informatica, CDAP • Databases: PostgreSQL/PostGIS, Oracle, MongoDB, Neo4j, MySQL'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Software Architecture & DevOps: Distributed Systems: Multi-tenant SaaS, microservices, API design, scalability • Geospatial Applications: Spatial algorithms, boundary estimation, clustering analysis • Web Applications: Full-stack development, RESTful APIs, real-time collaboration • DevOps: Docker, Vagrant, CI/CD (GitLab, GitHub), Celery, Airflow, nginx • Integration: Twilio API, WMS tile servers, CRM/DMP integration, OAuth'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'PROFESSIONAL EXPERIENCE'; Style = 'Heading 2'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'PARTNER & SENIOR SOFTWARE ENGINEER - Siege Analytics, Washington, DC | January 2014 – Present'; Style = 'Heading 3'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Geospatial Platform Architecture and Full-Stack Development'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Architected and engineered redistricting platform serving thousands of analysts with real-time collaborative editing, Census integration, and legal compliance analysis'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Developed boundary estimation microservice using incomplete data for boundary estimation without machine learning, processing geographies at national scale'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Built scalable ETL pipelines using PySpark and Sedona processing billions of geospatial records with sub-hour latency requirements'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Implemented advanced spatial clustering algorithms achieving 88% improvement in analytical targeting efficacy for political applications'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Created fraud detection systems processing multi-terabyte campaign finance datasets with real-time alerting capabilities'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Led technical architecture decisions integrating ESRI, OSGeo, and SAFE FME technologies for Fortune 500 and political clients'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'PRINCIPAL SOFTWARE ENGINEER - Clarity and Rigour, Washington, DC | 2012 – 2014'; Style = 'Heading 3'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Geospatial Solutions and Software Development'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Developed geospatial analysis frameworks and mapping applications for electoral research'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Built custom visualization tools and interactive dashboards for client presentations'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Architected data processing pipelines for large-scale demographic and geographic datasets'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Created web-based mapping applications using JavaScript, OpenLayers, and PostGIS'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'DIRECTOR OF DATA PRODUCTS - Helm, Washington, DC | 2010 – 2012'; Style = 'Heading 3'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Data Product Development and Engineering Leadership'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Led development of data-driven solutions and platform architecture for political organizations'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Managed cross-functional engineering teams building campaign management and voter targeting systems'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Designed and implemented scalable data platforms using Python, Django, and PostgreSQL'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Built RESTful APIs and microservices for campaign data integration'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'SENIOR SOFTWARE ENGINEER - GSD&M, Austin, TX | 2008 – 2010'; Style = 'Heading 3'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Campaign Technology and Analytics Development'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Developed voter targeting models and demographic analysis tools using Python and R'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Built web applications for campaign data visualization and reporting'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Created data integration systems connecting multiple campaign data sources'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Implemented machine learning algorithms for voter behavior prediction'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'SOFTWARE ENGINEER - Salsa Labs, Inc., Washington, DC | 2004 – 2006'; Style = 'Heading 3'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Political Technology Development'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Developed software solutions for political campaigns and advocacy groups using PHP and JavaScript'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Built web applications for voter engagement and campaign management'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Integrated third-party APIs and data sources for campaign tools'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Collaborated with political strategists to translate requirements into technical solutions'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 – 2004'; Style = 'Heading 3'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Nonprofit Technology Integration and Development'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Developed data management systems and web applications for social justice organizations'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Built custom applications for community engagement and advocacy using PHP and MySQL'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Provided technical training and support to nonprofit staff'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Integrated technology solutions within organizational frameworks'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'PROGRAMMER - Lake Research Partners, Washington, DC | 2001 – 2002'; Style = 'Heading 3'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Political Research and Data Analysis Tools'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Developed data analysis tools for political polling and research using Python and R'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Built statistical models and data visualization tools for research presentations'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Created automated reporting systems for survey data analysis'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Supported senior researchers with technical analysis and data processing'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 – 2001'; Style = 'Heading 3'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Political Field Operations and Data Management'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Developed data collection and management systems for field operations'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Built databases and reporting tools for campaign field work'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Trained field staff on data collection protocols and quality control'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Analyzed field data to inform campaign strategy and research findings'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'KEY ACHIEVEMENTS AND IMPACT'; Style = 'Heading 2'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = 'Geospatial Platform Development'; Style = 'Heading 3'; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Architected redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Built boundary estimation system achieving accurate geospatial results without machine learning using advanced PostGIS algorithms'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Developed econometric simulation platform with NetLogo multi-agent modeling and web interface'; Style = $null; Jc = $null; Bold = $false; Size = $null },
    @{ Text = '• Created comprehensive survey platform managing complete research lifecycle with integrated geospatial market segmentation'; Style = $null; Jc = $null; Bold = $false; Size = $null }
)

# --- Step 3: write out plain, unformatted text for every paragraph first ---
# (character formatting is applied afterwards so it cannot leak into later
#  paragraphs via Word's "inherit from previous run" behaviour)
$p = $d.Paragraphs.Item(1)
$p.Range.Text = $items[0].Text
for ($i = 1; $i -lt $items.Count; $i++) {
    $p.Range.InsertParagraphAfter()
    $p = $d.Paragraphs.Item($i + 1)
    $p.Range.InsertAfter($items[$i].Text)
}

# --- Step 4: second pass - apply paragraph styles / alignment / run formatting ---
$pos = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $it = $items[$i - 1]
    $txtLen = $it.Text.Length
    if ($it.Style) {
        $pp.Style = $it.Style
    }
    if ($it.Jc -eq "center") {
        $pp.Alignment = 1
    }
    if ($it.Bold -or $it.Size) {
        $tr = $d.Range($pos, $pos + $txtLen)
        if ($it.Bold) { $tr.Font.Bold = $true }
        if ($it.Size) { $tr.Font.Size = $it.Size }
    }
    $pos = $pp.Range.End
}

# --- Step 5: page margins (top/bottom 1440 twips = 72pt, left/right 1800 twips = 90pt) ---
$d.PageSetup.TopMargin = 72
$d.PageSetup.RightMargin = 90
$d.PageSetup.BottomMargin = 72
$d.PageSetup.LeftMargin = 90

Write-Output "Paragraphs written: $($d.Paragraphs.Count)"
